$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MCT-2A-M. T. R. M."
$ws.Range("F2").Value = "MCT-3A-Elementos de Máquinas"

$ws.Range("D3").Value = "MCT-2A-M. T. R. M."
$ws.Range("F3").Value = "MCT-3A-Elementos de Máquinas"

$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "MEC-2B-Mec. Tec. Res. Mat"
$ws.Range("F4").Value = "MCT-3A-Elementos de Máquinas"

$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "MEC-2B-Mec. Tec. Res. Mat"
$ws.Range("F6").Value = "MCT-3A-Elementos de Máquinas"

$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "MEC-2B-Mec. Tec. Res. Mat"
$ws.Range("F7").Value = "-"

$ws.Range("B8").Value = "-"
$ws.Range("D8").Value = "-"

$ws.Range("F10").Value = "MEC-2A-Mec. Tec. Res. Mat"

$ws.Range("B11").Value = "-"
$ws.Range("F11").Value = "MEC-2A-Mec. Tec. Res. Mat"

$ws.Range("B12").Value = "-"
$ws.Range("F12").Value = "MEC-2A-Mec. Tec. Res. Mat"

$ws.Range("B14").Value = "-"

$wb.Save()
